# The "zkratka" (abbreviation) column for a handful of rows had been
# mixed up between two rows sharing the same "nazev" (course name) /
# seminariciUcitIdno (teacher id) pair. Swap them back so each row's
# abbreviation matches its intended row (per commit: "Opravene kontroly
# hodin bez validnich ucitelu").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Podnikove informacni systemy: KEIS/EIS pairs
$ws.Range("B2").Value  = "EIS"
$ws.Range("B3").Value  = "KEIS"
$ws.Range("B4").Value  = "EIS"
$ws.Range("B5").Value  = "KEIS"
$ws.Range("B6").Value  = "EIS"
$ws.Range("B7").Value  = "KEIS"
$ws.Range("B8").Value  = "EIS"
$ws.Range("B9").Value  = "KEIS"
$ws.Range("B14").Value = "EIS"
$ws.Range("B15").Value = "KEIS"

# Odborna prezentace: KOPRE/OPRE pair (row 19, "Introduction to MATLAB" /
# ITM, is untouched)
$ws.Range("B18").Value = "OPRE"
$ws.Range("B20").Value = "KOPRE"
